$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge "This to" + bookmark + "ol is another way ..." into a single run.
#    (Concatenated text is already correct: "This to" + "ol is..." = "This tool is...")
#    Replacing the whole sentence collapses the runs/bookmark that used to sit
#    in the middle of the sentence into one contiguous run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "This tool is another way to measure WPI.  Wrap the yarn around the ruler section and the number of wraps is your WPI.  There are common names for various WPI measurements listed too. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "This tool is another way to measure WPI.  Wrap the yarn around the ruler section and the number of wraps is your WPI.  There are common names for various WPI measurements listed too. ",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Remove the entire "Bobbin Driver" tool section: the heading (with its
#    picture), the description paragraph and all five bulleted steps.
# ---------------------------------------------------------------------------
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13,[char]7) -eq "Bobbin Driver") {
        $headingPara = $i
        break
    }
}

$lastStepPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Load the bobbin by guiding the yarn")) {
        $lastStepPara = $i
    }
}

$rangeStart = $d.Paragraphs.Item($headingPara).Range.Start
$rangeEnd = $d.Paragraphs.Item($lastStepPara).Range.End
$d.Range($rangeStart, $rangeEnd).Delete()

# ---------------------------------------------------------------------------
# 3) Insert a new blank paragraph right after the "Spin Key" description
#    paragraph (the one that now ends in "...measurements listed too. ").
# ---------------------------------------------------------------------------
$spinKeyDescIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("This tool is another way to measure WPI")) {
        $spinKeyDescIndex = $i
        break
    }
}
$d.Paragraphs.Item($spinKeyDescIndex).Range.InsertParagraphAfter() | Out-Null

# ---------------------------------------------------------------------------
# 4) Fix up the final paragraph: drop the "Check out " lead-in and rework the
#    trailing text after the ElectricEelWheel.com hyperlink.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Check out ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
$d.Content.Find.Execute(" for more information and videos.", $false, $false, $false, $false, $false, $true, 1, $false, " has for more information and videos.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Append a brand-new paragraph pointing people at the Ravelry group.
# ---------------------------------------------------------------------------
$finalParaIndex = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($finalParaIndex)
$finalPara.Range.InsertParagraphAfter() | Out-Null

$newParaIndex = $finalParaIndex + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newRange = $newPara.Range
$newRange.Collapse(0)
$newRange.InsertAfter("www.ravelry.com/groups/electric-eel-wheel is a great online community for the EEW.")

# Turn the URL portion into a real hyperlink.
$hyperlinkStart = $newPara.Range.Start
$hyperlinkEnd = $hyperlinkStart + [string]"www.ravelry.com/groups/electric-eel-wheel".Length
$hyperlinkRange = $d.Range($hyperlinkStart, $hyperlinkEnd)
$d.Hyperlinks.Add($hyperlinkRange, "http://www.ravelry.com/groups/electric-eel-wheel", $null, $null, "www.ravelry.com/groups/electric-eel-wheel") | Out-Null

# ---------------------------------------------------------------------------
# 6) Relocate the (hidden) "_GoBack" bookmark to the very end of the document,
#    matching where Word leaves it after the last edit.
# ---------------------------------------------------------------------------
$endOfDoc = $d.Content.End - 1
$gobackRange = $d.Range($endOfDoc, $endOfDoc)
$d.Bookmarks.Add("_GoBack", $gobackRange) | Out-Null

Write-Output "done"
